$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped
# from 45206 to 45208 (i.e. +2 days) for every data row (rows 2-387).
$range = $ws.Range("C2:C387")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    $cell.Value2 = $cell.Value2 + 2
}
